# Fruta / hortaliza, semanal
# Insert a new weekly price row at row 504 (shifting existing rows 504-539
# down to 505-540), populated with the same field values as the row above
# (row 503) except for a new date (2023-12-05 -> serial 45265).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 504; default Insert() shifts the
# rows below it downward, exactly like the rows-inserted-below-2023-10-23
# row pattern visible in the rest of the sheet.
$ws.Rows.Item(504).Insert()

$ws.Cells.Item(504, 1).Value = 10
$ws.Cells.Item(504, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(504, 3).Value = "La Araucanía"
$ws.Cells.Item(504, 4).Value = 45265
$ws.Cells.Item(504, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(504, 5).Value = 9
$ws.Cells.Item(504, 6).Value = "Fruta"
$ws.Cells.Item(504, 7).Value = 100102
$ws.Cells.Item(504, 8).Value = "Cítricos"
$ws.Cells.Item(504, 9).Value = 100102006
$ws.Cells.Item(504, 10).Value = "Pomelo"
$ws.Cells.Item(504, 11).Value = "Start Ruby"
$ws.Cells.Item(504, 12).Value = "Primera"
$ws.Cells.Item(504, 13).Value = 50
$ws.Cells.Item(504, 14).Value = 14000
$ws.Cells.Item(504, 15).Value = 14000
$ws.Cells.Item(504, 16).Value = 14000
$ws.Cells.Item(504, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(504, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(504, 19).Value = 933
$ws.Cells.Item(504, 20).Value = 15
